$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" on every sheet that shows it
# (Overview!E2/F2 = zh-cn/de-de status, and the Status column on each
# per-language sheet).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# The shorter replacement text makes the status columns narrower; match
# the resized column widths left behind by the report regeneration.
$wsOverview.Range("E:F").ColumnWidth = 12.576851254417766
$wsZhCn.Range("C:C").ColumnWidth = 12.576851254417766
$wsDeDe.Range("C:C").ColumnWidth = 12.576851254417766
